# Agrego botón para simulación por forma reciente
# -> Agrega una columna "Resultado" (S) que indica si gano el Local,
#    la Visita o hubo Empate, en base a los goles de Local (D) y Visita (E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Calcular el resultado de cada partido (filas 2 a 147) comparando los
# goles de Local (columna D) contra los de Visita (columna E).
for ($i = 2; $i -le 147; $i++) {
    $golesLocal = $ws.Cells.Item($i, 4).Value()
    $golesVisita = $ws.Cells.Item($i, 5).Value()

    if ($golesLocal -gt $golesVisita) {
        $ws.Cells.Item($i, 19).Value = "Local"
    } elseif ($golesLocal -lt $golesVisita) {
        $ws.Cells.Item($i, 19).Value = "Visita"
    } else {
        $ws.Cells.Item($i, 19).Value = "Empate"
    }
}

# Encabezado de la nueva columna S, con el mismo estilo que el resto
# de encabezados (copiado desde R1).
$ws.Range("S1").Value = "Resultado"
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)

# Actualizar la vista activa de la hoja (desplazamiento y selección).
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("V11").Select()
